{"js": "// The document contains three adjacent runs that together read:\n//   <id>   (Courier New, color 7F6000, sz 9pt)\n//   p057v_2   (color 000000)\n//   </id>  (Courier New, color 7F6000, sz 9pt)\n// They need to be merged into a single run with text \"<id>p057v_2</id>\"\n// using the formatting of the first (\"<id>\") run.\n//\n// Searching for the full concatenated text returns a Range that spans all\n// three runs; replacing that range's text with itself collapses it back to\n// one run and keeps the formatting found at the start of the range (i.e.\n// the first run's formatting), which matches the target edit exactly.\n\nconst body = context.document.body;\nconst results = body.search(\"<id>p057v_2</id>\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text '<id>p057v_2</id>' in document body\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\"<id>p057v_2</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document contains three adjacent runs that together read:\n#   <id>      (Courier New, color 7F6000, sz 9pt)\n#   p057v_2   (color 000000)\n#   </id>     (Courier New, color 7F6000, sz 9pt)\n# They need to be merged into a single run with text \"<id>p057v_2</id>\"\n# using the formatting of the first (\"<id>\") run.\n#\n# Strategy: locate the first run's range (\"<id>\") and the full combined\n# range (\"<id>p057v_2</id>\"). Delete everything in the full range that comes\n# after the first run (i.e. \"p057v_2</id>\"), then re-append that same text\n# via InsertAfter on the still-intact first-run range. Word merges the\n# appended text into that original run, which keeps its formatting (and its\n# xml:space=\"preserve\"/rsid bookkeeping) instead of creating a brand new run.\n\n$d = $word.ActiveDocument\n\n$idRange = $d.Content\n$idFind = $idRange.Find\n$idFind.Text = \"<id>\"\n$idFound = $idFind.Execute()\n\nif (-not $idFound) {\n    throw \"Could not find '<id>' run in document\"\n}\n\n$fullRange = $d.Content\n$fullFind = $fullRange.Find\n$fullFind.Text = \"<id>p057v_2</id>\"\n$fullFound = $fullFind.Execute()\n\nif (-not $fullFound) {\n    throw \"Could not find '<id>p057v_2</id>' in document\"\n}\n\n$remainder = $fullRange.Text.Substring($idRange.End - $fullRange.Start)\n\n$restRange = $d.Range($idRange.End, $fullRange.End)\n$restRange.Delete()\n\n$idRange.InsertAfter($remainder)\n"}
